$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.13
$ws.Range("O2").Value = 1.69
$ws.Range("T2").Value = 1.05
$ws.Range("M3").Value = 1.11
$ws.Range("O3").Value = 1.63
$ws.Range("T3").Value = 1.08
$ws.Range("K4").Value = 1.8
$ws.Range("Q4").Value = 2.87
$ws.Range("R4").Value = 1.37
$ws.Range("G5").Value = 2.9
$ws.Range("I5").Value = 2.7
$ws.Range("K5").Value = 1.69
$ws.Range("L5").Value = 3.75
$ws.Range("Q5").Value = 3.6
$ws.Range("R5").Value = 1.25
$ws.Range("S5").Value = 9
$ws.Range("T5").Value = 1.07
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.98
$ws.Range("G6").Value = 2.15
$ws.Range("H6").Value = 2.82
$ws.Range("R6").Value = 1.3
$ws.Range("G7").Value = 2.05
$ws.Range("Q7").Value = 2.87
$ws.Range("R7").Value = 1.37
$ws.Range("G8").Value = 2.25
$ws.Range("H8").Value = 2.7
$ws.Range("R8").Value = 1.3
$ws.Range("G9").Value = 2.4
$ws.Range("R9").Value = 1.33
$ws.Range("S9").Value = 6.5
$ws.Range("T9").Value = 1.11
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2.65
$ws.Range("K10").Value = 1.77
$ws.Range("W10").Value = 2.37
$ws.Range("X10").Value = 1.5
$ws.Range("G11").Value = 1.96
$ws.Range("H11").Value = 2.75
$ws.Range("J11").Value = 2.87
$ws.Range("K11").Value = 1.69
$ws.Range("X11").Value = 1.33
$ws.Range("G12").Value = 2.15
$ws.Range("H12").Value = 2.82
$ws.Range("K12").Value = 1.77
$ws.Range("M12").Value = 1.14
$ws.Range("N12").Value = 5.5
$ws.Range("X12").Value = 1.47
$ws.Range("G13").Value = 2.5
$ws.Range("J13").Value = 3.4
$ws.Range("Z13").Value = 10
$ws.Range("AB13").Value = 23
$ws.Range("AI13").Value = 7
$ws.Range("AP15").Value = 1.74
$ws.Range("AQ15").Value = 1.99
$ws.Range("R27").Value = 1.47
$ws.Range("G28").Value = 2.62
$ws.Range("I28").Value = 2.37
$ws.Range("N28").Value = 8
$ws.Range("X28").Value = 1.58
$ws.Range("AC28").Value = 29
$ws.Range("AE28").Value = 7
$ws.Range("AL28").Value = 23
$ws.Range("AM28").Value = 23
$ws.Range("G29").Value = 2.75
$ws.Range("W29").Value = 1.58
$ws.Range("AJ29").Value = 13
$ws.Range("AM29").Value = 19
$ws.Range("AN29").Value = 26
$ws.Range("G30").Value = 1.81
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 2.5
$ws.Range("Q30").Value = 1.67
$ws.Range("R30").Value = 2.15
$ws.Range("Z30").Value = 10
$ws.Range("AB30").Value = 17
$ws.Range("AC30").Value = 15
$ws.Range("AI30").Value = 13
$ws.Range("AJ30").Value = 21
$ws.Range("G31").Value = 1.86
$ws.Range("H31").Value = 3
$ws.Range("I31").Value = 4.75
$ws.Range("J31").Value = 2.75
$ws.Range("K31").Value = 1.87
$ws.Range("Z31").Value = 7.5
$ws.Range("AJ31").Value = 21
$ws.Range("G32").Value = 1.96
$ws.Range("H32").Value = 3.2
$ws.Range("I32").Value = 3.9
$ws.Range("J32").Value = 2.75
$ws.Range("L32").Value = 4.33
$ws.Range("M32").Value = 1.07
$ws.Range("N32").Value = 9
$ws.Range("Q32").Value = 2.05
$ws.Range("R32").Value = 1.75
$ws.Range("S32").Value = 3.5
$ws.Range("T32").Value = 1.29
$ws.Range("Z32").Value = 9.5
$ws.Range("AA32").Value = 9
$ws.Range("AB32").Value = 17
$ws.Range("AC32").Value = 17
$ws.Range("AF32").Value = 6
$ws.Range("AI32").Value = 11
$ws.Range("AJ32").Value = 19
$ws.Range("AK32").Value = 13
$ws.Range("AL32").Value = 41
$ws.Range("AM32").Value = 34
$ws.Range("AO32").Value = 251
$ws.Range("G33").Value = 1.36
$ws.Range("H33").Value = 4.5
$ws.Range("I33").Value = 8.5
$ws.Range("J33").Value = 1.91
$ws.Range("L33").Value = 8
$ws.Range("M33").Value = 1.03
$ws.Range("N33").Value = 11
$ws.Range("O33").Value = 1.22
$ws.Range("Q33").Value = 1.85
$ws.Range("R33").Value = 2
$ws.Range("S33").Value = 3
$ws.Range("W33").Value = 2.1
$ws.Range("X33").Value = 1.67
$ws.Range("Y33").Value = 6.5
$ws.Range("AB33").Value = 8.5
$ws.Range("AC33").Value = 12
$ws.Range("AD33").Value = 29
$ws.Range("AF33").Value = 8.5
$ws.Range("AG33").Value = 21
$ws.Range("AH33").Value = 67
$ws.Range("AI33").Value = 19
$ws.Range("AK33").Value = 26
$ws.Range("AL33").Value = 101
$ws.Range("M34").Value = 1.08
$ws.Range("N34").Value = 6.5
$ws.Range("O34").Value = 1.5
$ws.Range("P34").Value = 2.37
$ws.Range("T34").Value = 1.11
$ws.Range("G35").Value = 2.5
$ws.Range("H35").Value = 3
$ws.Range("I35").Value = 2.9
$ws.Range("J35").Value = 3.1
$ws.Range("Q35").Value = 1.98
$ws.Range("R35").Value = 1.88
$ws.Range("S35").Value = 3.25
$ws.Range("T35").Value = 1.33
$ws.Range("U35").Value = 1.4
$ws.Range("V35").Value = 2.75
$ws.Range("Z35").Value = 12
$ws.Range("AB35").Value = 23
$ws.Range("N37").Value = 8
$ws.Range("O38").Value = 1.29
$ws.Range("P38").Value = 3.5
$ws.Range("Q38").Value = 1.95
$ws.Range("R38").Value = 1.9
